$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the Version value
$ws.Range("B3").Value = "0.7.0"

# Remove the "Jurisdiction" / "Chile" row (row 11), shifting rows 12-15 up
$ws.Rows(11).Delete()
